$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace "Dakota Myers" with "Franz Ferdinand" in B3
$ws.Range("B3").Value = "Franz Ferdinand"

# Update the selection to match the newly edited cell
$ws.Range("B3").Select()
